# Fix for BDF event import (edftype/type) idiosyncrasy.
# ICA_components.xlsx: the "B2" component-count value is no longer known up
# front (was a stale placeholder), so it is cleared out, and the shifted
# event-code bin total in "C2" is corrected from 5 to 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stale value in B2 but keep explicit (Normal) cell formatting so
# the cell continues to exist in the sheet (just blank) rather than being
# dropped entirely.
$ws.Cells.Item(2, 2).ClearContents()
$ws.Cells.Item(2, 2).Style = "Normal"

# Correct the shifted event-code bin total.
$ws.Cells.Item(2, 3).Value = 6

# Reflect where the user was last working in the sheet when they saved.
[void]$ws.Range("C8").Select()
